$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "TestCases"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TestCases"

# Add a new sheet "TestData" after TestCases
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TestData"

# ---- TestCases sheet content ----
$ws1.Range("A1").Value = "TestCases"
$ws1.Range("B1").Value = "Runmode"
$ws1.Range("A2").Value = "AddCustomerTest"
$ws1.Range("B2").Value = "Y"
$ws1.Range("A3").Value = "OpenAccountTest"
$ws1.Range("B3").Value = "N"

# COM ColumnWidth uses "character" units that get rounded to whole pixels
# internally (offset of 5/6 + nearest 1/6 step) by this engine, so feed it
# pre-compensated values to land on the OOXML widths from the target file
# (17 and ~9.57 -> closest reachable pixel-rounded value is 9.5).
$ws1.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws1.Columns.Item(2).ColumnWidth = 8.65

$ws1.Activate()
$ws1.Range("A3").Select()

# ---- TestData sheet content ----
$ws2.Range("A1").Value = "AddCustomerTest"

$ws2.Range("A2").Value = "Runmode"
$ws2.Range("B2").Value = "firstname"
$ws2.Range("C2").Value = "lastname"
$ws2.Range("D2").Value = "postcode"

$ws2.Range("A3").Value = "Y"
$ws2.Range("B3").Value = "manish"
$ws2.Range("C3").Value = "k"
$ws2.Range("D3").Value = "P6767"

$ws2.Range("A4").Value = "N"
$ws2.Range("B4").Value = "jyoti"
$ws2.Range("C4").Value = "k"
$ws2.Range("D4").Value = "X7878"

$ws2.Range("A6").Value = "OpenAccountTest"

$ws2.Range("A7").Value = "Runmode"
$ws2.Range("B7").Value = "customer"
$ws2.Range("C7").Value = "currency"

$ws2.Range("A8").Value = "Y"
$ws2.Range("B8").Value = "manish k"
$ws2.Range("C8").Value = "Rupee"

$ws2.Range("A9").Value = "Y"
$ws2.Range("B9").Value = "jyoti k"
$ws2.Range("C9").Value = "Dollar"

$ws2.Range("A6").Select()

$ws2.Activate()
